# Trade #65 closed at 2026-02-17 12:53:35 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet -------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value  = 1200.16   # Current Capital
$summary.Range("B4").Value  = 0.15      # Total P&L $
$summary.Range("B5").Value  = 0.05      # Total P&L %
$summary.Range("B6").Value  = 65        # Total Trades
$summary.Range("B7").Value  = 29        # Winning Trades
$summary.Range("B9").Value  = 44.62     # Win Rate %

# --- Strategy Status sheet (MarketMaking row) -----------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.16      # Capital
$status.Range("D4").Value = 65          # Trades
$status.Range("E4").Value = 0.15        # P&L $
$status.Range("F4").Value = 0.16        # P&L %
$status.Range("G4").Value = 44.62       # Win Rate %

# --- New trade row (#65) appended to "All Trades" and "MarketMaking" ------
$newRow = @{
    A = 65
    B = "2026-02-17"
    C = "12:53:29"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.87
    G = 0.92
    H = "CLOSED"
    I = 5.7471
    J = 0.05
    K = 100.16
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A66").Value = $newRow.A

    # The "Date" column holds a plain text value like "2026-02-17" in every
    # other row (not a real Excel date). Force text formatting so the
    # assignment isn't auto-converted into a date serial number, then
    # restore the default "Normal" style so the cell doesn't end up with a
    # stray explicit style index.
    $ws.Range("B66").NumberFormat = "@"
    $ws.Range("B66").Value = $newRow.B
    $ws.Range("B66").Style = "Normal"

    $ws.Range("C66").Value = $newRow.C
    $ws.Range("D66").Value = $newRow.D
    $ws.Range("E66").Value = $newRow.E
    $ws.Range("F66").Value = $newRow.F
    $ws.Range("G66").Value = $newRow.G
    $ws.Range("H66").Value = $newRow.H
    $ws.Range("I66").Value = $newRow.I
    $ws.Range("J66").Value = $newRow.J
    $ws.Range("K66").Value = $newRow.K
    $ws.Range("L66").Value = $newRow.L
    $ws.Range("M66").Value = $newRow.M
    $ws.Range("N66").Value = $newRow.N
    $ws.Range("O66").Value = $newRow.O
    $ws.Range("P66").Value = $newRow.P
    $ws.Range("Q66").Value = $newRow.Q
}
